# Insert a new snapshot column (CF) that holds the 2026-01-31 12:21:21
# price check, pushing the existing "nom" (CF->CG) and "url_produit"
# (CG->CH) columns one slot to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column 84 = CF. Inserting here shifts CF/CG -> CG/CH and copies the
# formatting of the column to the left, exactly like Excel's
# "Insert Column" command.
$ws.Columns.Item(84).Insert()

# Row 1 header: new timestamp for this snapshot.
$ws.Cells.Item(1, 84).Value = "2026-01-31 12:21:21"

# Rows 2-80 still carry a real last-seen price in column CE (column
# 83 - untouched by the insert at column 84). Copy that price forward
# into the new CF column so the history line continues flat until the
# next scrape records an actual new reading.
$lastPricedRow = 80
$srcPrices = $ws.Range($ws.Cells.Item(2, 83), $ws.Cells.Item($lastPricedRow, 83))
$dstPrices = $ws.Range($ws.Cells.Item(2, 84), $ws.Cells.Item($lastPricedRow, 84))
$dstPrices.Value2 = $srcPrices.Value2

# Rows 81-206 have no price history left (CE is blank there already),
# so the newly inserted CF cells for those rows stay blank too - no
# action needed, Insert() already left them empty.
